$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 72 (pushes existing rows 72-194 down to 73-195)
$ws.Rows.Item(72).Insert()

# Populate the new row 72 with the new weekly price record.
$ws.Range("A72").Value = 5
$ws.Range("B72").Value = "Macroferia Regional de Talca"
$ws.Range("C72").Value = "Maule"
$ws.Range("D72").Value = 44477
$ws.Range("E72").Value = 7
$ws.Range("F72").Value = 100114013
$ws.Range("G72").Value = "Zanahoria"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 500
$ws.Range("K72").Value = 7000
$ws.Range("L72").Value = 7000
$ws.Range("M72").Value = 7000
$ws.Range("N72").Value = "$/saco 20 kilos"
$ws.Range("O72").Value = "Región de Ñuble"
$ws.Range("P72").Value = 350
$ws.Range("Q72").Value = 20
$ws.Range("R72").Value = "Hortaliza"
